# Update NATMI LR-pair sheet (Ccl12-Ackr2) with new TPM-based results.
# - Row 2 (Inflammatory-Mac): refresh TPM-derived metrics
# - Row 3: sending cluster changes from Neutrophils -> MuSCs, refresh metrics
# - Row 4: sending cluster becomes Neutrophils (shifted down), refresh metrics
# - Row 5 (new): Resolving-Mac -> Ccl12/Ackr2/FAPs interaction added

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Inflammatory-Mac ----
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.70362466666668
$ws.Range("H2").Value = 128.110874
$ws.Range("I2").Value = 0.510021191154308
$ws.Range("J2").Value = 0.5102913077099245
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2786473333333334
$ws.Range("N2").Value = 0.8359420000000001
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 11.89925113703423
$ws.Range("R2").Value = 107.093260233308
$ws.Range("S2").Value = 0.510021191154308
$ws.Range("T2").Value = 0.5102913077099245

# ---- Row 3: MuSCs (previously Neutrophils) ----
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.132963
$ws.Range("H3").Value = 0.265926
$ws.Range("I3").Value = 0.001588013855235666
$ws.Range("J3").Value = 0.001059236597621443
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2786473333333334
$ws.Range("N3").Value = 0.8359420000000001
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.037049785382
$ws.Range("R3").Value = 0.222298712292
$ws.Range("S3").Value = 0.001588013855235666
$ws.Range("T3").Value = 0.001059236597621443

# ---- Row 4: Neutrophils ----
$ws.Range("A4").Value = "Neutrophils"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.925701999999999
$ws.Range("H4").Value = 26.777106
$ws.Range("I4").Value = 0.1066021257320059
$ws.Range("J4").Value = 0.1066585841684857
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2786473333333334
$ws.Range("N4").Value = 0.8359420000000001
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 2.487123060428
$ws.Range("R4").Value = 22.384107543852
$ws.Range("S4").Value = 0.1066021257320059
$ws.Range("T4").Value = 0.1066585841684857

# ---- Row 5: Resolving-Mac (new row) ----
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 31.96682866666667
$ws.Range("H5").Value = 95.900486
$ws.Range("I5").Value = 0.3817886692584505
$ws.Range("J5").Value = 0.3819908715239683
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2786473333333334
$ws.Range("N5").Value = 0.8359420000000001
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 8.907471563090223
$ws.Range("R5").Value = 80.16724406781201
$ws.Range("S5").Value = 0.3817886692584505
$ws.Range("T5").Value = 0.3819908715239683
